$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '38.318.06'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.24%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.095.42'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.94%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '228.69'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.20%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '61.07'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.06%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.379'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.18%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0849'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.13%  '
$ws.Range("E11").Value = '  +0.05%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.407.94'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.03%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.73'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.48%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '22.30'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +5.72%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.47'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +5.39%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.776'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.02%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.088.71'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.80%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '38.205.98'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.09%  '
$ws.Range("B19").Value = 'Litecoin'
$ws.Range("C19").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '70.36'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.78%  '
$ws.Range("B20").Value = 'Uniswap'
$ws.Range("C20").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.00'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.23%  '
$ws.Range("E21").Value = '  +1.12%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '224.39'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.23%  '
$ws.Range("E23").Value = '  +0.01%  '
$ws.Range("E24").Value = '  +0.59%  '
$ws.Range("E25").Value = '  +2.18%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '169.91'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.46%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.41'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.08%  '
$ws.Range("E28").Value = '  +0.49%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.98'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.67%  '
$ws.Range("E31").Value = '  -0.76%  '
$ws.Range("E32").Value = '  +8.48%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.73'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.96%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.43'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.49%  '
$ws.Range("E35").Value = '  -0.36%  '
$ws.Range("E36").Value = '  +4.19%  '
$ws.Range("E37").Value = '  +1.25%  '
$ws.Range("E38").Value = '  +5.80%  '
$ws.Range("E39").Value = '  -0.10%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.13'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.78%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.547.08'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.70%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '99.82'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.71%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0218'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.66%  '
$ws.Range("E44").Value = '  +1.05%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0907'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.91%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.15'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.28%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.11'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.68%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.51'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.82%  '
$ws.Range("E49").Value = '  +1.85%  '
$ws.Range("E50").Value = '  +0.94%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.294.19'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.02%  '
